$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.763.33"
$ws.Range("E2").Value = "  -3.09%  "

$ws.Range("D3").Value = "2.577.89"
$ws.Range("E3").Value = "  -5.21%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.17%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +1.21%  "

$ws.Range("E9").Value = "  -3.11%  "

$ws.Range("E10").Value = "  -1.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.32%  "

$ws.Range("E12").Value = "  -2.51%  "

$ws.Range("D13").Value = "3.030.81"
$ws.Range("E13").Value = "  -5.33%  "

$ws.Range("E14").Value = "  -4.79%  "

$ws.Range("D15").Value = "61.667.10"
$ws.Range("E15").Value = "  -2.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000143"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.30%  "

$ws.Range("D17").Value = "2.575.70"
$ws.Range("E17").Value = "  -5.33%  "

$ws.Range("E18").Value = "  -5.51%  "

$ws.Range("E19").Value = "  -2.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "337.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.00%  "

$ws.Range("E21").Value = "  -5.97%  "

$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.493"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.48%  "

$ws.Range("E25").Value = "  -1.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  -1.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.93%  "

$ws.Range("E29").Value = "  -4.60%  "

$ws.Range("E30").Value = "  -2.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.73%  "

$ws.Range("E34").Value = "  -2.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.36%  "

$ws.Range("E36").Value = "  -4.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "333.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.928"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.71%  "

$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").Value = "2.128.26"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("E46").Value = "  -3.75%  "

$ws.Range("E47").Value = "  -1.08%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.22%  "

$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0545"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.40%  "

$ws.Range("E50").Value = "  -1.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0238"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.77%  "
